{"js": "// Adds a \"Stuff to do:\" subtitle followed by three new bullet items\n// (\"Logika pro battle\", \"Grafika\", \"Hudba\") at the end of the document,\n// reusing the pre-existing trailing empty list paragraph for the first\n// bullet (same as the authored diff).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document ends with an empty \"List Paragraph\" (Odstavecseseznamem)\n// item; it becomes the \"Logika pro battle\" bullet.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// New subtitle paragraph right before that trailing bullet.\nconst heading = lastParagraph.insertParagraph(\"Stuff to do:\", \"Before\");\nheading.style = \"Subtitle\";\n\n// Fill the previously-empty trailing bullet paragraph.\nlastParagraph.insertText(\"Logika pro battle\", \"Replace\");\n\n// Two more bullet items after it, joining the same numbered list\n// (numId 1, level 0) as the rest of the bullets.\nconst grafika = lastParagraph.insertParagraph(\"Grafika\", \"After\");\ngrafika.style = \"List Paragraph\";\ngrafika.attachToList(1, 0);\n\nconst hudba = grafika.insertParagraph(\"Hudba\", \"After\");\nhudba.style = \"List Paragraph\";\nhudba.attachToList(1, 0);\n\nawait context.sync();\n", "ps1": "# Adds a \"Stuff to do:\" subtitle followed by three new bullet items\n# (\"Logika pro battle\", \"Grafika\", \"Hudba\") at the end of the document,\n# reusing the pre-existing trailing empty list paragraph for the first\n# bullet (same as the authored diff).\n\n$d = $word.ActiveDocument\n\n# The document ends with an empty \"List Paragraph\" (Odstavecseseznamem)\n# bullet; split a new empty paragraph in right before it (inherits the\n# same List Paragraph style/numbering).\n$lastPara = $d.Paragraphs.Last\n$lastPara.Range.InsertParagraphBefore()\n\n# That freshly inserted paragraph is now the second-to-last one; turn it\n# into the \"Stuff to do:\" subtitle.\n$count = $d.Paragraphs.Count\n$heading = $d.Paragraphs.Item($count - 1)\n$heading.Style = \"Subtitle\"\n$heading.Range.Text = \"Stuff to do:\"\n\n# Fill the original trailing (now last) bullet paragraph.\n$battle = $d.Paragraphs.Item($count)\n$battle.Range.Text = \"Logika pro battle\"\n\n# Add \"Grafika\" bullet after it.\n$battle.Range.InsertParagraphAfter()\n$grafika = $d.Paragraphs.Item($count + 1)\n$grafika.Range.Text = \"Grafika\"\n\n# Add \"Hudba\" bullet after that.\n$grafika.Range.InsertParagraphAfter()\n$hudba = $d.Paragraphs.Item($count + 2)\n$hudba.Range.Text = \"Hudba\"\n"}
